$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.514.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.52%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.844.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.07%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -1.06%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'334.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.39%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.10%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4617"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.62%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3866"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.87%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'45.94"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.92%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.07906"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.18%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.9990"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.64%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'21.51"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.06%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'5.968"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.00%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.846.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.18%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'7.148"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.12%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.009"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.18%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'88.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.94%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.06669"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.29%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.00001034"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.32%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'17.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.77%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'1.007"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.14%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'27.512.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.62%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.392"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.00%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'10.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.55%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.305"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.91%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.070.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.00%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'159.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.24%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'19.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.79%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.114"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +3.03%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'5.410"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.23%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'120.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.42%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.9746"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.11%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.09398"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.56%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'3.601"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.78%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'5.298"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.03%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.334"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.69%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.06012"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.11%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.02222"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.09%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'8.277"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.65%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.181"
$ws.Range("D40").Style = "Normal"
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = "'0.5904"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.00%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = "'0.1863"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.41%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = "'10.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.39%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = 'WEMIXTOKEN'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").Value = "'1.240"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.40%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = "'0.5584"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.52%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'12.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.17%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = "'1.909"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.34%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = "'0.06698"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.83%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = "'110.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.59%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").Value = "'1.048"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.04%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = 'PaxDollar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D51").Value = "'1.007"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.18%  "
$ws.Range("E51").Style = "Normal"
